# Generate Report for Handoff
# The localization-status report moves from "In Translation" to
# "Ready for handoff": status cells are updated, the handoff timestamps are
# refreshed, and the Status column is widened to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refresh the "latest handoff xliff generate" timestamps -----------------
$overview.Range("G2").Value = "2016-08-21 02:44:16"
$zhcn.Range("H2").Value     = "2016-08-21 02:44:12"
$dede.Range("H2").Value     = "2016-08-21 02:44:16"

# --- Widen the Status columns so the longer text isn't truncated ------------
$overview.Range("E1").ColumnWidth = 16.333333333333332
$overview.Range("F1").ColumnWidth = 16.333333333333332
$zhcn.Range("C1").ColumnWidth     = 16.333333333333332
$dede.Range("C1").ColumnWidth     = 16.333333333333332
